$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry data which differs between rows 8, 9 and 10
# (A, B, D, E, F, G, H, Q, R, Z, AB). Columns Y/AA (dates) and all other
# columns are identical across these three rows, so they are left untouched.
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 26, 28)

function Get-RowValues($row) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value()
    }
    return $vals
}

# Capture the "before" contents of the three affected rows first, since we
# will be overwriting them in place.
$row8 = Get-RowValues 8
$row9 = Get-RowValues 9
$row10 = Get-RowValues 10

function Set-RowValues($row, $vals) {
    foreach ($c in $cols) {
        $v = $vals[$c]
        if ($v -eq $null) { $v = "" }
        $ws.Cells.Item($row, $c).Value = $v
    }
}

# Cyclic shift: old row8 -> row9, old row9 -> row10, old row10 -> row8
Set-RowValues 8 $row10
Set-RowValues 9 $row8
Set-RowValues 10 $row9
